$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3495
$ws.Range("I62").Value = 2993.3333
$ws.Range("J62").Value = 5000
$ws.Range("K62").Value = 2993.3333
$ws.Range("L62").Value = 5000
$ws.Range("M62").Value = -2369.3333
$ws.Range("N62").Value = -6248

$ws.Range("H65").Value = 3495
$ws.Range("I65").Value = 2993.3333
$ws.Range("J65").Value = 5000
$ws.Range("K65").Value = 14966.6665
$ws.Range("L65").Value = 25000
$ws.Range("M65").Value = -11846.6665
$ws.Range("N65").Value = -31240

$ws.Range("H69").Value = 3447.5

$ws.Range("H72").Value = 3447.5

$ws.Range("H96").Value = 863.44446
$ws.Range("I96").Value = 382.35715
$ws.Range("J96").Value = 2547.25
$ws.Range("K96").Value = 1147.07145
$ws.Range("L96").Value = 7641.75
$ws.Range("M96").Value = 225.9285500000001
$ws.Range("N96").Value = -10387.75

$ws.Range("H101").Value = 636.4545000000001
$ws.Range("I101").Value = 402.2857
$ws.Range("J101").Value = 1046.25
$ws.Range("K101").Value = 1206.8571
$ws.Range("L101").Value = 3138.75
$ws.Range("M101").Value = 415.1428999999998
$ws.Range("N101").Value = -6382.75

$ws.Range("H137").Value = 1376.5555
$ws.Range("I137").Value = 1348.5333
$ws.Range("J137").Value = 1516.6666
$ws.Range("K137").Value = 4045.5999
$ws.Range("L137").Value = 4549.9998
$ws.Range("M137").Value = -1495.5999
$ws.Range("N137").Value = -9649.9998

$ws.Range("H138").Value = 2395.3542
$ws.Range("I138").Value = 1347.0454
$ws.Range("J138").Value = 3282.3845
$ws.Range("K138").Value = 4041.1362
$ws.Range("L138").Value = 9847.1535
$ws.Range("M138").Value = 1098.8638
$ws.Range("N138").Value = -20127.1535

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 24121.35
$ws.Range("I2").Value = 992
$ws.Range("J2").Value = 91406.73
$ws.Range("K2").Value = 992
$ws.Range("L2").Value = 91406.73
$ws.Range("M2").Value = -879
$ws.Range("N2").Value = -91632.73

$ws.Range("H61").Value = 1466.1774
$ws.Range("I61").Value = 1008.3469
$ws.Range("K61").Value = 1008.3469
$ws.Range("M61").Value = -796.3469

$ws.Range("H74").Value = 745.4358999999999
$ws.Range("I74").Value = 696.625
$ws.Range("J74").Value = 968.5714
$ws.Range("K74").Value = 696.625
$ws.Range("L74").Value = 968.5714
$ws.Range("M74").Value = 177.375
$ws.Range("N74").Value = -2716.5714

$ws.Range("H77").Value = 745.4358999999999
$ws.Range("I77").Value = 696.625
$ws.Range("J77").Value = 968.5714
$ws.Range("K77").Value = 3483.125
$ws.Range("L77").Value = 4842.857
$ws.Range("M77").Value = 884.875
$ws.Range("N77").Value = -13578.857

$ws.Range("H116").Value = 24121.35
$ws.Range("I116").Value = 992
$ws.Range("J116").Value = 91406.73
$ws.Range("K116").Value = 992
$ws.Range("L116").Value = 91406.73
$ws.Range("M116").Value = 1302
$ws.Range("N116").Value = -95994.73

$ws.Range("H132").Value = 16779.553
$ws.Range("I132").Value = 18859.875
$ws.Range("J132").Value = 5684.5
$ws.Range("K132").Value = 56579.625
$ws.Range("L132").Value = 17053.5
$ws.Range("M132").Value = -54049.625
$ws.Range("N132").Value = -22113.5

$ws.Range("H136").Value = 1466.1774
$ws.Range("I136").Value = 1008.3469
$ws.Range("K136").Value = 3025.0407
$ws.Range("M136").Value = -475.0407

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 24121.35
$ws.Range("I3").Value = 992
$ws.Range("J3").Value = 91406.73
$ws.Range("K3").Value = 992
$ws.Range("L3").Value = 91406.73
$ws.Range("M3").Value = -878
$ws.Range("N3").Value = -91634.73

$ws.Range("H105").Value = 97023.336
$ws.Range("I105").Value = 68513.92999999999
$ws.Range("J105").Value = 168296.83
$ws.Range("K105").Value = 68513.92999999999
$ws.Range("L105").Value = 168296.83
$ws.Range("M105").Value = -66766.92999999999
$ws.Range("N105").Value = -171790.83

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 26835.316
$ws.Range("I31").Value = 499
$ws.Range("J31").Value = 52263.484
$ws.Range("K31").Value = 499
$ws.Range("L31").Value = 52263.484
$ws.Range("M31").Value = -204
$ws.Range("N31").Value = -52853.484

$ws.Range("H34").Value = 26835.316
$ws.Range("I34").Value = 499
$ws.Range("J34").Value = 52263.484
$ws.Range("K34").Value = 499
$ws.Range("L34").Value = 52263.484
$ws.Range("M34").Value = -297
$ws.Range("N34").Value = -52667.484

$ws.Range("H58").Value = 2535.1875
$ws.Range("I58").Value = 2422.6365
$ws.Range("K58").Value = 2422.6365
$ws.Range("M58").Value = -2219.6365

$ws.Range("H99").Value = 8030.647
$ws.Range("I99").Value = 1762
$ws.Range("K99").Value = 1762
$ws.Range("M99").Value = -264

$ws.Range("H126").Value = 8030.647
$ws.Range("I126").Value = 1762
$ws.Range("K126").Value = 5286
$ws.Range("M126").Value = -2816

$ws.Range("H134").Value = 1836.1333
$ws.Range("I134").Value = 1644.5
$ws.Range("J134").Value = 2123.5833
$ws.Range("K134").Value = 4933.5
$ws.Range("L134").Value = 6370.749899999999
$ws.Range("M134").Value = -2398.5
$ws.Range("N134").Value = -11440.7499

$ws.Range("H136").Value = 2535.1875
$ws.Range("I136").Value = 2422.6365
$ws.Range("K136").Value = 7267.9095
$ws.Range("M136").Value = -4717.9095

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1611.4445
$ws.Range("I75").Value = 917.1667
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 2751.5001
$ws.Range("L75").Value = 9000
$ws.Range("M75").Value = -1753.5001
$ws.Range("N75").Value = -10996

$ws.Range("H78").Value = 1611.4445
$ws.Range("I78").Value = 917.1667
$ws.Range("J78").Value = 3000
$ws.Range("K78").Value = 8254.5003
$ws.Range("L78").Value = 27000
$ws.Range("M78").Value = -3262.5003
$ws.Range("N78").Value = -36984

$ws.Range("H113").Value = 1015.2381
$ws.Range("J113").Value = 563.8461
$ws.Range("L113").Value = 1691.5383
$ws.Range("N113").Value = -6031.5383

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H76").Value = 9999.5
$ws.Range("J76").Value = 9999.5
$ws.Range("L76").Value = 9999.5
$ws.Range("N76").Value = -10675.5

$ws.Range("H79").Value = 9999.5
$ws.Range("J79").Value = 9999.5
$ws.Range("L79").Value = 9999.5
$ws.Range("N79").Value = -12339.5

$ws.Range("H98").Value = 29999.5
$ws.Range("J98").Value = 29999.5
$ws.Range("L98").Value = 29999.5
$ws.Range("N98").Value = -35989.5

$ws.Range("H106").Value = 30275
$ws.Range("J106").Value = 30275
$ws.Range("L106").Value = 30275
$ws.Range("N106").Value = -32799
